$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B70 was stored as an inline string "4"; change it to a real number 4.
$ws.Range("B70").Value = 4

# Append a new row 71 with the additional annotation data.
$ws.Range("A71").Value = "Ying Tang"
$ws.Range("B71").Value = "'3"
$ws.Range("B71").Style = "Normal"
$ws.Range("C71").Value = "We will include,for clarity"
$ws.Range("D71").Value = "SUG"
$ws.Range("E71").Value = "WRI"
$ws.Range("F71").Value = "3c70bee3-3ebe-492b-b68a-cb43e1a99f35"
$ws.Range("G71").Value = "H1Ww66x0-_annotated.xlsx"
$ws.Range("H71").Value = "We will include additional details on the hyper-parameters of the baselines for clarity."
